$wb = $excel.ActiveWorkbook

# --- Insert a new "Player Info" sheet before "ODI Batting" ---
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingForInsert)
$playerInfo.Name = "Player Info"

# Re-resolve sheet references by name now that the collection has shifted
# (references captured before Add() point at positions, not identities).
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# Header row (bold / centered / bordered, matching the other sheets' header style)
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").Value = "'4737"
$playerInfo.Range("B2").Value = "Mahedi Hasan"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE (column D) ---
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4452"
$batting.Range("D3").Value = "'4453"
$batting.Range("D4").Value = "'4455"

# --- Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE (column B) ---
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4452"
$bowling.Range("B3").Value = "'4453"
$bowling.Range("B4").Value = "'4455"
